$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 33.00563866666667
$ws.Range("H2").Value = 99.01691600000001
$ws.Range("I2").Value = 0.9169150302490913
$ws.Range("J2").Value = 0.9169150302490912
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.386787
$ws.Range("N2").Value = 7.160361
$ws.Range("O2").Value = 0.2248318673684458
$ws.Range("P2").Value = 0.2248318673684458
$ws.Range("Q2").Value = 78.77742929629734
$ws.Range("R2").Value = 708.9968636666761
$ws.Range("S2").Value = 0.2061517184690982
$ws.Range("T2").Value = 0.2061517184690982

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 33.00563866666667
$ws.Range("H3").Value = 99.01691600000001
$ws.Range("I3").Value = 0.9169150302490913
$ws.Range("J3").Value = 0.9169150302490912
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 4.875491666666666
$ws.Range("N3").Value = 14.626475
$ws.Range("O3").Value = 0.4592642308492391
$ws.Range("P3").Value = 0.4592642308492391
$ws.Range("Q3").Value = 160.9187162723445
$ws.Range("R3").Value = 1448.2684464511
$ws.Range("S3").Value = 0.4211062761214557
$ws.Range("T3").Value = 0.4211062761214557

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 33.00563866666667
$ws.Range("H4").Value = 99.01691600000001
$ws.Range("I4").Value = 0.9169150302490913
$ws.Range("J4").Value = 0.9169150302490912
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.8868146666666666
$ws.Range("N4").Value = 2.660444
$ws.Range("O4").Value = 0.08353665304712675
$ws.Range("P4").Value = 0.08353665304712675
$ws.Range("Q4").Value = 29.26988445230045
$ws.Range("R4").Value = 263.428960070704
$ws.Range("S4").Value = 0.07659601275561406
$ws.Range("T4").Value = 0.07659601275561406

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 33.00563866666667
$ws.Range("H5").Value = 99.01691600000001
$ws.Range("I5").Value = 0.9169150302490913
$ws.Range("J5").Value = 0.9169150302490912
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.466781666666666
$ws.Range("N5").Value = 7.400345
$ws.Range("O5").Value = 0.2323672487351883
$ws.Range("P5").Value = 0.2323672487351882
$ws.Range("Q5").Value = 81.41770435955777
$ws.Range("R5").Value = 732.7593392360201
$ws.Range("S5").Value = 0.2130610229029233
$ws.Range("T5").Value = 0.2130610229029232

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.07465466666666666
$ws.Range("H6").Value = 0.223964
$ws.Range("I6").Value = 0.002073948231580021
$ws.Range("J6").Value = 0.00207394823158002
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.386787
$ws.Range("N6").Value = 7.160361
$ws.Range("O6").Value = 0.2248318673684458
$ws.Range("P6").Value = 0.2248318673684458
$ws.Range("Q6").Value = 0.1781847878893333
$ws.Range("R6").Value = 1.603663091004
$ws.Range("S6").Value = 0.0004662896537316221
$ws.Range("T6").Value = 0.0004662896537316219

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.07465466666666666
$ws.Range("H7").Value = 0.223964
$ws.Range("I7").Value = 0.002073948231580021
$ws.Range("J7").Value = 0.00207394823158002
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 4.875491666666666
$ws.Range("N7").Value = 14.626475
$ws.Range("O7").Value = 0.4592642308492391
$ws.Range("P7").Value = 0.4592642308492391
$ws.Range("Q7").Value = 0.363978205211111
$ws.Range("R7").Value = 3.2758038469
$ws.Range("S7").Value = 0.0009524902393977379
$ws.Range("T7").Value = 0.0009524902393977377

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.07465466666666666
$ws.Range("H8").Value = 0.223964
$ws.Range("I8").Value = 0.002073948231580021
$ws.Range("J8").Value = 0.00207394823158002
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.8868146666666666
$ws.Range("N8").Value = 2.660444
$ws.Range("O8").Value = 0.08353665304712675
$ws.Range("P8").Value = 0.08353665304712675
$ws.Range("Q8").Value = 0.0662048533351111
$ws.Range("R8").Value = 0.595843680016
$ws.Range("S8").Value = 0.0001732506938592023
$ws.Range("T8").Value = 0.0001732506938592022

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.07465466666666666
$ws.Range("H9").Value = 0.223964
$ws.Range("I9").Value = 0.002073948231580021
$ws.Range("J9").Value = 0.00207394823158002
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.466781666666666
$ws.Range("N9").Value = 7.400345
$ws.Range("O9").Value = 0.2323672487351883
$ws.Range("P9").Value = 0.2323672487351882
$ws.Range("Q9").Value = 0.1841567630644444
$ws.Range("R9").Value = 1.65741086758
$ws.Range("S9").Value = 0.0004819176445914585
$ws.Range("T9").Value = 0.0004819176445914583

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.576299333333333
$ws.Range("H10").Value = 7.728898
$ws.Range("I10").Value = 0.07157103078692272
$ws.Range("J10").Value = 0.0715710307869227
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 2.386787
$ws.Range("N10").Value = 7.160361
$ws.Range("O10").Value = 0.2248318673684458
$ws.Range("P10").Value = 0.2248318673684458
$ws.Range("Q10").Value = 6.149077756908667
$ws.Range("R10").Value = 55.341699812178
$ws.Range("S10").Value = 0.01609144850130836
$ws.Range("T10").Value = 0.01609144850130836

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2.576299333333333
$ws.Range("H11").Value = 7.728898
$ws.Range("I11").Value = 0.07157103078692272
$ws.Range("J11").Value = 0.0715710307869227
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 4.875491666666666
$ws.Range("N11").Value = 14.626475
$ws.Range("O11").Value = 0.4592642308492391
$ws.Range("P11").Value = 0.4592642308492391
$ws.Range("Q11").Value = 12.56072593050555
$ws.Range("R11").Value = 113.04653337455
$ws.Range("S11").Value = 0.03287001440544327
$ws.Range("T11").Value = 0.03287001440544327

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 2.576299333333333
$ws.Range("H12").Value = 7.728898
$ws.Range("I12").Value = 0.07157103078692272
$ws.Range("J12").Value = 0.0715710307869227
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.8868146666666666
$ws.Range("N12").Value = 2.660444
$ws.Range("O12").Value = 0.08353665304712675
$ws.Range("P12").Value = 0.08353665304712675
$ws.Range("Q12").Value = 2.284700034523556
$ws.Range("R12").Value = 20.562300310712
$ws.Range("S12").Value = 0.00597880436707239
$ws.Range("T12").Value = 0.005978804367072389

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 2.576299333333333
$ws.Range("H13").Value = 7.728898
$ws.Range("I13").Value = 0.07157103078692272
$ws.Range("J13").Value = 0.0715710307869227
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 2.466781666666666
$ws.Range("N13").Value = 7.400345
$ws.Range("O13").Value = 0.2323672487351883
$ws.Range("P13").Value = 0.2323672487351882
$ws.Range("Q13").Value = 6.355167963312222
$ws.Range("R13").Value = 57.19651166981
$ws.Range("S13").Value = 0.01663076351309869
$ws.Range("T13").Value = 0.01663076351309868

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.3398056666666667
$ws.Range("H14").Value = 1.019417
$ws.Range("I14").Value = 0.009439990732406145
$ws.Range("J14").Value = 0.009439990732406145
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 2.386787
$ws.Range("N14").Value = 7.160361
$ws.Range("O14").Value = 0.2248318673684458
$ws.Range("P14").Value = 0.2248318673684458
$ws.Range("Q14").Value = 0.8110437477263334
$ws.Range("R14").Value = 7.299393729537
$ws.Range("S14").Value = 0.002122410744307697
$ws.Range("T14").Value = 0.002122410744307697

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.3398056666666667
$ws.Range("H15").Value = 1.019417
$ws.Range("I15").Value = 0.009439990732406145
$ws.Range("J15").Value = 0.009439990732406145
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 4.875491666666666
$ws.Range("N15").Value = 14.626475
$ws.Range("O15").Value = 0.4592642308492391
$ws.Range("P15").Value = 0.4592642308492391
$ws.Range("Q15").Value = 1.656719696119444
$ws.Range("R15").Value = 14.910477265075
$ws.Range("S15").Value = 0.004335450082942453
$ws.Range("T15").Value = 0.004335450082942453

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.3398056666666667
$ws.Range("H16").Value = 1.019417
$ws.Range("I16").Value = 0.009439990732406145
$ws.Range("J16").Value = 0.009439990732406145
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.8868146666666666
$ws.Range("N16").Value = 2.660444
$ws.Range("O16").Value = 0.08353665304712675
$ws.Range("P16").Value = 0.08353665304712675
$ws.Range("Q16").Value = 0.3013446490164444
$ws.Range("R16").Value = 2.712101841148
$ws.Range("S16").Value = 0.0007885852305811041
$ws.Range("T16").Value = 0.0007885852305811041

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 0.3398056666666667
$ws.Range("H17").Value = 1.019417
$ws.Range("I17").Value = 0.009439990732406145
$ws.Range("J17").Value = 0.009439990732406145
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 2.466781666666666
$ws.Range("N17").Value = 7.400345
$ws.Range("O17").Value = 0.2323672487351883
$ws.Range("P17").Value = 0.2323672487351882
$ws.Range("Q17").Value = 0.8382263887627777
$ws.Range("R17").Value = 7.544037498864999
$ws.Range("S17").Value = 0.002193544674574891
$ws.Range("T17").Value = 0.002193544674574891
